# Change default fuel tax revenue GRA setting.
#
# The "fuel tax revenue" weights live in row 9 (B9:F9) of the
# "Set Values Here" sheet. The "GRA-fueltax" sheet reads those same
# weights back in via an array formula (TRANSPOSE), so updating the
# source row is sufficient to ripple the new numbers through.

$wb = $excel.ActiveWorkbook

# --- Update the "fuel tax revenue" weights on "Set Values Here" ---
$setValuesSheet = $wb.Worksheets.Item("Set Values Here")
$setValuesSheet.Range("C9").Value = 5
$setValuesSheet.Range("D9").Value = 0
$setValuesSheet.Range("F9").Value = 5

# Leave the cursor where the author left it when done editing.
$setValuesSheet.Range("F10").Select()

# --- Touch the dependent "GRA-fueltax" sheet, leaving the cursor on B2 ---
$fuelTaxSheet = $wb.Worksheets.Item("GRA-fueltax")
$fuelTaxSheet.Range("B2").Select()

# --- Restore focus to the sheet that was active originally ---
$aboutSheet = $wb.Worksheets.Item("About")
$aboutSheet.Activate()
